$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 743.7143
$ws.Range("I2").Value = 34.333332
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 34.333332
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = 78.666668
$ws.Range("N2").Value = -5226
$ws.Range("H20").Value = 1132.8572
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("H35").Value = 1132.8572
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("H48").Value = 4612.5
$ws.Range("J48").Value = 5250
$ws.Range("L48").Value = 15750
$ws.Range("N48").Value = -16334
$ws.Range("H56").Value = 4612.5
$ws.Range("J56").Value = 5250
$ws.Range("L56").Value = 15750
$ws.Range("N56").Value = -16818
$ws.Range("H62").Value = 5527.1816
$ws.Range("I62").Value = 4942.857
$ws.Range("K62").Value = 4942.857
$ws.Range("M62").Value = -4318.857
$ws.Range("H64").Value = 1875
$ws.Range("I64").Value = 1875
$ws.Range("K64").Value = 1875
$ws.Range("M64").Value = -1627
$ws.Range("H65").Value = 5527.1816
$ws.Range("I65").Value = 4942.857
$ws.Range("K65").Value = 24714.285
$ws.Range("M65").Value = -21594.285
$ws.Range("H67").Value = 1875
$ws.Range("I67").Value = 1875
$ws.Range("K67").Value = 1875
$ws.Range("M67").Value = -1017
$ws.Range("H76").Value = 5282.2
$ws.Range("I76").Value = 4805.5
$ws.Range("K76").Value = 4805.5
$ws.Range("M76").Value = -4490.5
$ws.Range("H79").Value = 5282.2
$ws.Range("I79").Value = 4805.5
$ws.Range("K79").Value = 4805.5
$ws.Range("M79").Value = -3713.5
$ws.Range("H106").Value = 6832.6665
$ws.Range("I106").Value = 6832.6665
$ws.Range("K106").Value = 6832.6665
$ws.Range("M106").Value = -6201.6665
$ws.Range("H112").Value = 2359.375
$ws.Range("I112").Value = 1237.5
$ws.Range("J112").Value = 2733.3333
$ws.Range("K112").Value = 3712.5
$ws.Range("L112").Value = 8199.999899999999
$ws.Range("M112").Value = -2604.5
$ws.Range("N112").Value = -10415.9999
$ws.Range("H116").Value = 6909.3335
$ws.Range("I116").Value = 5897.5
$ws.Range("K116").Value = 5897.5
$ws.Range("M116").Value = -2455.5
$ws.Range("H132").Value = 4113.636
$ws.Range("I132").Value = 4147.4443
$ws.Range("K132").Value = 12442.3329
$ws.Range("M132").Value = -9912.332900000001
$ws.Range("H137").Value = 2052.3333
$ws.Range("I137").Value = 1912.1428
$ws.Range("K137").Value = 5736.428400000001
$ws.Range("M137").Value = -3186.428400000001
$ws.Range("N20").ClearContents()
$ws.Range("N35").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 17663.834
$ws.Range("J13").Value = 17663.834
$ws.Range("L13").Value = 17663.834
$ws.Range("N13").Value = -17951.834
$ws.Range("H41").Value = 14014
$ws.Range("I41").Value = 2018.6666
$ws.Range("J41").Value = 50000
$ws.Range("K41").Value = 2018.6666
$ws.Range("L41").Value = 50000
$ws.Range("M41").Value = -1604.6666
$ws.Range("N41").Value = -50828
$ws.Range("H132").Value = 904
$ws.Range("I132").Value = 904
$ws.Range("K132").Value = 2712
$ws.Range("M132").Value = -182
$ws.Range("H139").Value = 45107.5
$ws.Range("J139").Value = 45107.5
$ws.Range("L139").Value = 45107.5
$ws.Range("N139").Value = -55387.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1899
$ws.Range("I20").Value = 1199
$ws.Range("J20").Value = 2132.3333
$ws.Range("K20").Value = 1199
$ws.Range("L20").Value = 2132.3333
$ws.Range("M20").Value = -952
$ws.Range("N20").Value = -2626.3333
$ws.Range("H54").Value = 50000
$ws.Range("I54").Value = 50000
$ws.Range("K54").Value = 50000
$ws.Range("M54").Value = -49516
$ws.Range("H86").Value = 4050
$ws.Range("I86").Value = 3950
$ws.Range("J86").Value = 4150
$ws.Range("K86").Value = 3950
$ws.Range("L86").Value = 4150
$ws.Range("M86").Value = -2827
$ws.Range("N86").Value = -6396
$ws.Range("H89").Value = 4050
$ws.Range("I89").Value = 3950
$ws.Range("J89").Value = 4150
$ws.Range("K89").Value = 19750
$ws.Range("L89").Value = 20750
$ws.Range("M89").Value = -14134
$ws.Range("N89").Value = -31982
$ws.Range("H94").Value = 1969.6
$ws.Range("I94").Value = 1899.5
$ws.Range("K94").Value = 1899.5
$ws.Range("M94").Value = -1448.5
$ws.Range("H105").Value = 3060
$ws.Range("I105").Value = 2933.3333
$ws.Range("K105").Value = 2933.3333
$ws.Range("M105").Value = -1186.3333
$ws.Range("H107").Value = 1143
$ws.Range("I107").Value = 1143
$ws.Range("K107").Value = 1143
$ws.Range("M107").Value = 777
$ws.Range("H130").Value = 59999.668
$ws.Range("J130").Value = 59999.668
$ws.Range("L130").Value = 59999.668
$ws.Range("N130").Value = -70039.66800000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 30998.8
$ws.Range("J59").Value = 34373.125
$ws.Range("L59").Value = 34373.125
$ws.Range("N59").Value = -36663.125
$ws.Range("H109").Value = 50698.25
$ws.Range("J109").Value = 49844.668
$ws.Range("L109").Value = 49844.668
$ws.Range("N109").Value = -51924.668

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 150.8
$ws.Range("I2").Value = 112.375
$ws.Range("K2").Value = 674.25
$ws.Range("M2").Value = -561.25
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("H75").Value = 1143.3334
$ws.Range("I75").Value = 30
$ws.Range("J75").Value = 1700
$ws.Range("K75").Value = 90
$ws.Range("L75").Value = 5100
$ws.Range("M75").Value = 908
$ws.Range("N75").Value = -7096
$ws.Range("H78").Value = 1143.3334
$ws.Range("I78").Value = 30
$ws.Range("J78").Value = 1700
$ws.Range("K78").Value = 270
$ws.Range("L78").Value = 15300
$ws.Range("M78").Value = 4722
$ws.Range("N78").Value = -25284
$ws.Range("H140").Value = 1595.6666
$ws.Range("I140").Value = 1114.8
$ws.Range("J140").Value = 4000
$ws.Range("K140").Value = 3344.4
$ws.Range("L140").Value = 12000
$ws.Range("M140").Value = 1835.6
$ws.Range("N140").Value = -22360
$ws.Range("H141").Value = 2004.5
$ws.Range("I141").Value = 2004.5
$ws.Range("K141").Value = 6013.5
$ws.Range("M141").Value = -833.5
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 19003
$ws.Range("J20").Value = 19003
$ws.Range("L20").Value = 19003
$ws.Range("N20").Value = -19493
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("H69").Value = 92498.75
$ws.Range("J69").Value = 92498.75
$ws.Range("L69").Value = 92498.75
$ws.Range("N69").Value = -93996.75
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("H72").Value = 92498.75
$ws.Range("J72").Value = 92498.75
$ws.Range("L72").Value = 277496.25
$ws.Range("N72").Value = -284984.25
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("H80").Value = 3912.2222
$ws.Range("I80").Value = 2486.6667
$ws.Range("K80").Value = 2486.6667
$ws.Range("M80").Value = -1488.6667
$ws.Range("H83").Value = 3912.2222
$ws.Range("I83").Value = 2486.6667
$ws.Range("K83").Value = 12433.3335
$ws.Range("M83").Value = -7441.333500000001
$ws.Range("H102").Value = 1379.8334
$ws.Range("I102").Value = 1055.8
$ws.Range("K102").Value = 1055.8
$ws.Range("M102").Value = 566.2
$ws.Range("H122").Value = 31300504
$ws.Range("I122").Value = 41700336
$ws.Range("K122").Value = 125101008
$ws.Range("M122").Value = -125098558
$ws.Range("M24").ClearContents()
$ws.Range("M59").ClearContents()
$ws.Range("M70").ClearContents()
$ws.Range("M73").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3643.889
$ws.Range("J46").Value = 3999.375
$ws.Range("L46").Value = 3999.375
$ws.Range("N46").Value = -4375.375
$ws.Range("H123").Value = 77999.5
$ws.Range("J123").Value = 77999.5
$ws.Range("L123").Value = 77999.5
$ws.Range("N123").Value = -87799.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 50000
$ws.Range("J50").Value = 50000
$ws.Range("L50").Value = 50000
$ws.Range("N50").Value = -51262
$ws.Range("H96").Value = 1233.1666
$ws.Range("I96").Value = 1079.8
$ws.Range("K96").Value = 1079.8
$ws.Range("M96").Value = 293.2
$ws.Range("H109").Value = 52000
$ws.Range("J109").Value = 69000
$ws.Range("L109").Value = 69000
$ws.Range("N109").Value = -71774
$ws.Range("H122").Value = 768.8
$ws.Range("I122").Value = 768.8
$ws.Range("K122").Value = 2306.4
$ws.Range("M122").Value = 143.6000000000004
$ws.Range("H132").Value = 1999
$ws.Range("I132").Value = 1999
$ws.Range("K132").Value = 5997
$ws.Range("M132").Value = -3467
